# Make the "ORIGEN" field dynamic by adding it as its own column (F)
# instead of it living elsewhere - mirrors the existing header/data row
# layout (A:E) so the new column matches the look of the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the first header cell (bold white text on a
# black fill with a border) onto the new header cell, then enter the
# header text and the data value underneath it.
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "ORIGEN"

# The header style used across the table has a border on the top/left/
# right (xlEdgeTop/xlEdgeLeft/xlEdgeRight); the new trailing column only
# keeps the left/right edges, so drop the top/bottom rules explicitly.
$ws.Range("F1").Borders.Item(8).LineStyle = -4142   # xlEdgeTop -> xlLineStyleNone
$ws.Range("F1").Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> xlLineStyleNone

$ws.Range("F2").Value = "AGP AMBA"

$ws.Range("F2").Select()
